# Update F-column ("想去人数" / want-to-go count) figures to the
# latest scrape snapshot (gh-pages data refresh, commit 456a3b4).
# "展览"/"演出"/"本地生活" are the per-category sheets; "全部类型" is the
# combined sheet that mirrors the same rows, so it gets the same updates.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 113  # was 114
$ws.Range("F4").Value = 582  # was 575
$ws.Range("F6").Value = 25  # was 24
$ws.Range("F7").Value = 1923  # was 1918
$ws.Range("F8").Value = 5238  # was 5203
$ws.Range("F9").Value = 1460  # was 1453
$ws.Range("F11").Value = 3061  # was 3053
$ws.Range("F12").Value = 695  # was 694
$ws.Range("F13").Value = 33  # was 32
$ws.Range("F14").Value = 1275  # was 1268
$ws.Range("F15").Value = 4185  # was 4164
$ws.Range("F16").Value = 990  # was 985
$ws.Range("F18").Value = 1639  # was 1637
$ws.Range("F19").Value = 2585  # was 2580
$ws.Range("F20").Value = 26  # was 25
$ws.Range("F21").Value = 15  # was 14
$ws.Range("F22").Value = 121  # was 119
$ws.Range("F23").Value = 141  # was 140
$ws.Range("F24").Value = 955  # was 951
$ws.Range("F25").Value = 285  # was 283
$ws.Range("F26").Value = 76  # was 75
$ws.Range("F27").Value = 73  # was 72
$ws.Range("F28").Value = 200  # was 199
$ws.Range("F29").Value = 1064  # was 1062
$ws.Range("F30").Value = 349  # was 343
$ws.Range("F31").Value = 26  # was 22
$ws.Range("F32").Value = 114  # was 113
$ws.Range("F34").Value = 223  # was 212
$ws.Range("F35").Value = 1626  # was 1614
$ws.Range("F36").Value = 2146  # was 2141
$ws.Range("F37").Value = 1000  # was 997
$ws.Range("F40").Value = 592  # was 590
$ws.Range("F41").Value = 261  # was 257
$ws.Range("F43").Value = 646  # was 645
$ws.Range("F44").Value = 388  # was 386
$ws.Range("F45").Value = 294  # was 289
$ws.Range("F46").Value = 203  # was 202
$ws.Range("F47").Value = 127  # was 126

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 143  # was 142

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 721  # was 720

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 721  # was 720
$ws.Range("F6").Value = 582  # was 575
$ws.Range("F7").Value = 25  # was 24
$ws.Range("F8").Value = 1923  # was 1918
$ws.Range("F9").Value = 5238  # was 5203
$ws.Range("F10").Value = 1460  # was 1453
$ws.Range("F13").Value = 3061  # was 3053
$ws.Range("F14").Value = 33  # was 32
$ws.Range("F15").Value = 1275  # was 1268
$ws.Range("F16").Value = 4185  # was 4164
$ws.Range("F17").Value = 990  # was 985
$ws.Range("F18").Value = 1639  # was 1637
$ws.Range("F20").Value = 2585  # was 2580
$ws.Range("F22").Value = 26  # was 25
$ws.Range("F25").Value = 141  # was 140
$ws.Range("F26").Value = 143  # was 142
$ws.Range("F27").Value = 955  # was 951
$ws.Range("F28").Value = 285  # was 283
$ws.Range("F29").Value = 73  # was 72
$ws.Range("F30").Value = 200  # was 199
$ws.Range("F32").Value = 1064  # was 1062
$ws.Range("F33").Value = 349  # was 343
$ws.Range("F34").Value = 26  # was 22
$ws.Range("F36").Value = 1626  # was 1614
$ws.Range("F37").Value = 2146  # was 2141
$ws.Range("F38").Value = 1000  # was 997
$ws.Range("F42").Value = 592  # was 590
$ws.Range("F43").Value = 261  # was 257
$ws.Range("F44").Value = 646  # was 645
$ws.Range("F45").Value = 388  # was 386
$ws.Range("F46").Value = 294  # was 289
$ws.Range("F47").Value = 203  # was 202
$ws.Range("F48").Value = 127  # was 126
